$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the .jpg image-file extensions with .png across the stimuli list (A2:A16)
for ($r = 2; $r -le 16; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    $new = $old -replace '\.jpg$', '.png'
    $cell.Value = $new
}

# Update the active selection to match the author's final cursor position
$ws.Range("G16").Select()
